# Update cryptocurrency price/volume data per Mon Jun 10 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.374.07"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "3.673.43"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'639.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.58%  "
$ws.Range("D6").Value = "'159.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "'7.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").Value = "'0.449"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "4.294.71"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "'32.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "3.672.25"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "69.358.06"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "'15.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D20").Value = "'465.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "'9.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'0.647"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "'79.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").Value = "3.820.05"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("D27").Value = "'10.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").Value = "'9.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").Value = "'2.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").Value = "'1.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("D31").Value = "'2.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("D35").Value = "'6.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("D36").Value = "3.666.54"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").Value = "'8.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'5.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.63%  "
$ws.Range("D40").Value = "'177.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.39%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'0.0898"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").Value = "'0.925"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").Value = "'2.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").Value = "'27.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.18%  "
$ws.Range("D48").Value = "'0.000268"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("E49").Value = "  -3.71%  "
$ws.Range("D50").Value = "'7.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("E51").Value = "  -3.77%  "
